$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Break the link between C3 and A1: hardcode the formula to 70 instead of =A1
$ws.Range("C3").Formula = "=70"

# Force a full recalculation so dependent cells (C4:C38, E, K, C40, C41, C43...) update
$excel.CalculateFullRebuild()

# Update the view: scroll so row 10 is at the top and select D42
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D42").Select()
